# Insert a new data row at row 248 (pushes the old rows 248..351 down to 249..352)
# and populate it with the new "Early Burlat" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(248).Insert()

$ws.Cells.Item(248, 1).Value2  = 6
$ws.Cells.Item(248, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(248, 3).Value2  = "Metropolitana"
$ws.Cells.Item(248, 4).Value2  = 44510
$ws.Cells.Item(248, 5).Value2  = 13
$ws.Cells.Item(248, 6).Value2  = "Fruta"
$ws.Cells.Item(248, 7).Value2  = 100103
$ws.Cells.Item(248, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(248, 9).Value2  = 100103001
$ws.Cells.Item(248, 10).Value2 = "Cereza"
$ws.Cells.Item(248, 11).Value2 = "Early Burlat"
$ws.Cells.Item(248, 12).Value2 = "Primera"
$ws.Cells.Item(248, 13).Value2 = 96
$ws.Cells.Item(248, 14).Value2 = 39000
$ws.Cells.Item(248, 15).Value2 = 42000
$ws.Cells.Item(248, 16).Value2 = 40500
$ws.Cells.Item(248, 17).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(248, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(248, 19).Value2 = 2700
$ws.Cells.Item(248, 20).Value2 = 15
